$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("ACzFALBrQjSkcqxWdGRx6x4eGZbf2RbwMtC7nS5ckBVP", 90, "asefuaihwefcbebn"),
    @("ACzFALBrQjSkcqxWdGRx6x4eGZbf2RbwMtC7nS5ckBVP", 90, "asefuaihwefcbebn"),
    @("ACzFALBrQjSkcqxWdGRx6x4eGZbf2RbwMtC7nS5ckBVP", 90, "asefuaihwefcbebn"),
    @("ACzFALBrQjSkcqxWdGRx6x4eGZbf2RbwMtC7nS5ckBVP", 90, "asefuaihwefcbebn"),
    @("ACzFALBrQjSkcqxWdGRx6x4eGZbf2RbwMtC7nS5ckBVP", 90, "56ZcbkWFcsBagDXNB8Yyk6nyCdTa6S8i71wVz3zdwVwvtgmAAJh1cVxyLx6P3FkwSRqawQHmQNgM2iitv3RoEM4"),
    @("ACzFALBrQjSkcqxWdGRx6x4eGZbf2RbwMtC7nS5ckBVP", 90, "hHPtn2TXQMt1YcbMKTYLMf8hmR6NFekF3jyjRrMqNWr48ULoCG9cFeyTA98C2rBWcEf4kMuTW42tVoNw5eN7FHJ")
)

$r = 7
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
